$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.212.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.481.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.54%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.520'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0810'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.871.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.508.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.839'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.107.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0931'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.46%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.89%  '
$ws.Range("E29").Value = '  +0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.08%  '
$ws.Range("E31").Value = '  +5.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("E34").Value = '  +3.42%  '
$ws.Range("E35").Value = '  +2.62%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.23%  '
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.62%  '
$ws.Range("E40").Value = '  +1.71%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '120.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("E44").Value = '  +2.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.948.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +13.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.59%  '
